# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff-handback timestamp
# cells on the Overview, zh-cn and de-de sheets to reflect a freshly
# generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-10-18 12:16:15"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-10-18 12:16:03"
$zhcn.Range("K3").Value = "2016-10-18 12:16:59"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-10-18 12:16:15"
$dede.Range("K3").Value = "2016-10-18 12:17:17"
